# PHOENIX-5876: Completed the creation of the new water connection
#
# On the "approvalDetails" sheet, duplicate the "commissioner1" approval row
# so it also appears right after the "commissioner" row (in addition to its
# original position further down), and record the officer name for the
# "deputyExecutiveEngineer" approval row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("approvalDetails")

# Insert a new row above row 6 (the "engineer" row), pushing everything
# from row 6 onward down by one.
$ws.Rows.Item(6).Insert()

# Populate the newly inserted row 6 with the same data as the
# "commissioner1" approval row (now at row 13).
$ws.Range("A6").Value = "commissioner1"
$ws.Range("B6").Value = "ADMINISTRATION"
$ws.Range("C6").Value = "Commissioner"
$ws.Range("D6").Value = "S.Ravindra Babu/ADM_Commissioner_1"
$ws.Range("E6").Value = "Forward to commissioner"

# Record the officer name for the "deputyExecutiveEngineer" row (now row 14).
$ws.Range("D14").Value = "S.Nayab Rasool/ENG_Dy. Executive Engineer_1"

$ws.Range("D18").Select()
